# Insert a new weekly price-report row above the current row 503 for
# "Feria Lagunitas de Puerto Montt - Betarraga" (Hortaliza), pushing the
# existing rows 503:551 down to 504:552.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 503:551 down to 504:552, leaving row 503 free for the new entry.
$ws.Rows.Item(503).Insert()

# Fill in the new row 503 with this week's data.
$ws.Range("A503").Value = 4
$ws.Range("B503").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C503").Value = "Los Lagos"
$ws.Range("D503").Value = 45223
$ws.Range("E503").Value = 10
$ws.Range("F503").Value = 100114014
$ws.Range("G503").Value = "Betarraga"
$ws.Range("H503").Value = "Sin especificar"
$ws.Range("I503").Value = "Primera"
$ws.Range("J503").Value = 1000
$ws.Range("K503").Value = 1000
$ws.Range("L503").Value = 1100
$ws.Range("M503").Value = 1050
$ws.Range("N503").Value = "`$/paquete 5 unidades"
$ws.Range("O503").Value = "Región Metropolitana"
$ws.Range("P503").Value = 210
$ws.Range("Q503").Value = 5
$ws.Range("R503").Value = "Hortaliza"
